$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 3938.7778
$ws.Range("J29").Value = 5706.5
$ws.Range("L29").Value = 17119.5
$ws.Range("N29").Value = -17681.5

$ws.Range("H31").Value = 5471.4
$ws.Range("I31").Value = 4089.5
$ws.Range("K31").Value = 12268.5
$ws.Range("M31").Value = -12038.5

$ws.Range("H76").Value = 4650
$ws.Range("I76").Value = 3950
$ws.Range("J76").Value = 5350
$ws.Range("K76").Value = 3950
$ws.Range("L76").Value = 5350
$ws.Range("M76").Value = -3635
$ws.Range("N76").Value = -5980

$ws.Range("H79").Value = 4650
$ws.Range("I79").Value = 3950
$ws.Range("J79").Value = 5350
$ws.Range("K79").Value = 3950
$ws.Range("L79").Value = 5350
$ws.Range("M79").Value = -2858
$ws.Range("N79").Value = -7534

$ws.Range("H86").Value = 7545.6665
$ws.Range("I86").Value = 9909.299999999999
$ws.Range("J86").Value = 2818.4
$ws.Range("K86").Value = 9909.299999999999
$ws.Range("L86").Value = 2818.4
$ws.Range("M86").Value = -8786.299999999999
$ws.Range("N86").Value = -5064.4

$ws.Range("H87").Value = 35994.8
$ws.Range("J87").Value = 39993.75
$ws.Range("L87").Value = 39993.75
$ws.Range("N87").Value = -42489.75

$ws.Range("H89").Value = 7545.6665
$ws.Range("I89").Value = 9909.299999999999
$ws.Range("J89").Value = 2818.4
$ws.Range("K89").Value = 49546.5
$ws.Range("L89").Value = 14092
$ws.Range("M89").Value = -43930.5
$ws.Range("N89").Value = -25324

$ws.Range("H90").Value = 35994.8
$ws.Range("J90").Value = 39993.75
$ws.Range("L90").Value = 119981.25
$ws.Range("N90").Value = -132461.25

$ws.Range("H137").Value = 4829.3477
$ws.Range("I137").Value = 4780.6924
$ws.Range("K137").Value = 14342.0772
$ws.Range("M137").Value = -11792.0772

$ws.Range("H138").Value = 7296.0737
$ws.Range("I138").Value = 4472.485
$ws.Range("J138").Value = 8343.022000000001
$ws.Range("K138").Value = 13417.455
$ws.Range("L138").Value = 25029.066
$ws.Range("M138").Value = -8277.454999999998
$ws.Range("N138").Value = -35309.06600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1933.6842
$ws.Range("I32").Value = 1932.5
$ws.Range("K32").Value = 1932.5
$ws.Range("M32").Value = -1645.5

$ws.Range("H45").Value = 1349.8334
$ws.Range("I45").Value = 1259.8
$ws.Range("K45").Value = 1259.8
$ws.Range("M45").Value = -882.8

$ws.Range("H61").Value = 66669480
$ws.Range("I61").Value = 71431304
$ws.Range("K61").Value = 71431304
$ws.Range("M61").Value = -71431092

$ws.Range("H132").Value = 3962069.8
$ws.Range("I132").Value = 2567573
$ws.Range("K132").Value = 7702719
$ws.Range("M132").Value = -7700189

$ws.Range("H136").Value = 66669480
$ws.Range("I136").Value = 71431304
$ws.Range("K136").Value = 214293912
$ws.Range("M136").Value = -214291362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1762.8422
$ws.Range("I20").Value = 1620.5834
$ws.Range("K20").Value = 1620.5834
$ws.Range("M20").Value = -1373.5834

$ws.Range("H76").Value = 22500
$ws.Range("J76").Value = 20000
$ws.Range("L76").Value = 20000
$ws.Range("N76").Value = -20630

$ws.Range("H79").Value = 22500
$ws.Range("J79").Value = 20000
$ws.Range("L79").Value = 20000
$ws.Range("N79").Value = -22184

$ws.Range("H134").Value = 15613440
$ws.Range("I134").Value = 15613440
$ws.Range("K134").Value = 46840320
$ws.Range("M134").Value = -46837785

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 130.27272
$ws.Range("I7").Value = 149.25
$ws.Range("K7").Value = 149.25
$ws.Range("M7").Value = -36.25

$ws.Range("H9").Value = 27400
$ws.Range("J9").Value = 27400
$ws.Range("L9").Value = 27400
$ws.Range("N9").Value = -27736

$ws.Range("H22").Value = 3316.1516
$ws.Range("J22").Value = 101.5
$ws.Range("L22").Value = 101.5
$ws.Range("N22").Value = -801.5

$ws.Range("H31").Value = 3640.3035
$ws.Range("I31").Value = 2416.1277
$ws.Range("K31").Value = 2416.1277
$ws.Range("M31").Value = -2121.1277

$ws.Range("H34").Value = 3640.3035
$ws.Range("I34").Value = 2416.1277
$ws.Range("K34").Value = 2416.1277
$ws.Range("M34").Value = -2214.1277

$ws.Range("H107").Value = 2967766
$ws.Range("I107").Value = 3623688.2
$ws.Range("K107").Value = 3623688.2
$ws.Range("M107").Value = -3621768.2

$ws.Range("H132").Value = 38463284
$ws.Range("J132").Value = 2049
$ws.Range("L132").Value = 6147
$ws.Range("N132").Value = -11207

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 924781.4399999999
$ws.Range("I4").Value = 1200990.9
$ws.Range("K4").Value = 3602972.7
$ws.Range("M4").Value = -3602860.7

$ws.Range("H86").Value = 531
$ws.Range("I86").Value = 455.91666
$ws.Range("J86").Value = 606.0833
$ws.Range("K86").Value = 1367.74998
$ws.Range("L86").Value = 1818.2499
$ws.Range("M86").Value = -181.7499800000001
$ws.Range("N86").Value = -4190.2499

$ws.Range("H89").Value = 531
$ws.Range("I89").Value = 455.91666
$ws.Range("J89").Value = 606.0833
$ws.Range("K89").Value = 4103.24994
$ws.Range("L89").Value = 5454.7497
$ws.Range("M89").Value = 1824.75006
$ws.Range("N89").Value = -17310.7497

$ws.Range("H107").Value = 1260.5
$ws.Range("I107").Value = 671
$ws.Range("J107").Value = 1555.25
$ws.Range("K107").Value = 2013
$ws.Range("L107").Value = 4665.75
$ws.Range("M107").Value = -93
$ws.Range("N107").Value = -8505.75

$ws.Range("H113").Value = 167029.67
$ws.Range("J113").Value = 447.25
$ws.Range("L113").Value = 1341.75
$ws.Range("N113").Value = -5681.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2652.2222
$ws.Range("J80").Value = 2395
$ws.Range("L80").Value = 2395
$ws.Range("N80").Value = -4391

$ws.Range("H83").Value = 2652.2222
$ws.Range("J83").Value = 2395
$ws.Range("L83").Value = 11975
$ws.Range("N83").Value = -21959

$ws.Range("H113").Value = 35061.242
$ws.Range("I113").Value = 41635.6
$ws.Range("K113").Value = 41635.6
$ws.Range("M113").Value = -39465.6

$ws.Range("H122").Value = 5106.4614
$ws.Range("I122").Value = 3046.1052
$ws.Range("K122").Value = 9138.3156
$ws.Range("M122").Value = -6688.3156

$ws.Range("H132").Value = 2672851.5
$ws.Range("I132").Value = 2853841
$ws.Range("K132").Value = 8561523
$ws.Range("M132").Value = -8558993

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 22680
$ws.Range("I45").Value = 9020.5
$ws.Range("J45").Value = 49999
$ws.Range("K45").Value = 9020.5
$ws.Range("L45").Value = 49999
$ws.Range("M45").Value = -8613.5
$ws.Range("N45").Value = -50813

$ws.Range("H82").Value = 2842.6667
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 2842.6667
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 2842.6667
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -3564.6667

$ws.Range("H85").Value = 2842.6667
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 2842.6667
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 2842.6667
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -5338.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 3266.3333
$ws.Range("I100").Value = 4399.5
$ws.Range("K100").Value = 8799
$ws.Range("M100").Value = -8258

$ws.Range("H126").Value = 1110
$ws.Range("I126").Value = 1113.3334
$ws.Range("J126").Value = 1100
$ws.Range("K126").Value = 3340.0002
$ws.Range("L126").Value = 3300
$ws.Range("M126").Value = -870.0001999999999
$ws.Range("N126").Value = -8240

$ws.Range("H136").Value = 27780122
$ws.Range("I136").Value = 29414040
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 88242120
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -88239570
$ws.Range("N136").Value = -15600
